# Weekly data refresh: two new price observations were inserted at the top
# of this week's block (rows 169-170), pushing all subsequent rows down by
# two and growing the used range from A1:R245 to A1:R247.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 169 (existing rows 169.. shift down to 171..)
$ws.Rows.Item(169).Insert()
$ws.Rows.Item(169).Insert()

# --- New row 169 ---
$ws.Range("A169").Value = 9
$ws.Range("B169").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C169").Value = "Metropolitana"
$ws.Range("D169").Value = 44455
$ws.Range("E169").Value = 13
$ws.Range("F169").Value = 100112031
$ws.Range("G169").Value = "Poroto verde"
$ws.Range("H169").Value = "Magnum"
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 38
$ws.Range("K169").Value = 33000
$ws.Range("L169").Value = 34000
$ws.Range("M169").Value = 33500
$ws.Range("N169").Value = "$/malla 25 kilos"
$ws.Range("O169").Value = "Perú"
$ws.Range("P169").Value = 1340
$ws.Range("Q169").Value = 25
$ws.Range("R169").Value = "Hortaliza"

# --- New row 170 ---
$ws.Range("A170").Value = 9
$ws.Range("B170").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C170").Value = "Metropolitana"
$ws.Range("D170").Value = 44455
$ws.Range("E170").Value = 13
$ws.Range("F170").Value = 100112031
$ws.Range("G170").Value = "Poroto verde"
$ws.Range("H170").Value = "Sin especificar"
$ws.Range("I170").Value = "Primera"
$ws.Range("J170").Value = 18
$ws.Range("K170").Value = 32000
$ws.Range("L170").Value = 33000
$ws.Range("M170").Value = 32500
$ws.Range("N170").Value = "$/malla 25 kilos"
$ws.Range("O170").Value = "Perú"
$ws.Range("P170").Value = 1300
$ws.Range("Q170").Value = 25
$ws.Range("R170").Value = "Hortaliza"

Write-Output "done"
